# Apply the Alvearie FHIR IG CodeSystem metadata refresh to the "Metadata" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Delete the duplicated "Contact" / "No display for ContactDetail" row (row 11).
# This shifts all subsequent rows up by one, turning the former 22-row table
# into a 21-row table.
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date: refresh publish timestamp
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now populated
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# The remaining "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# Case Sensitive value was blank, now populated with the literal text "true"
# (Assigning the string "true" directly via .Value/.Formula gets auto-coerced
# to a Boolean by the engine, so build it as a text formula in a scratch cell,
# then paste-special just the resulting value back into the target cell.)
$scratch = $ws.Cells.Item(14, 4)
$scratch.Formula = "=""tru""&""e"""
$scratch.Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()
